$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Map of row -> new value for column F (dSF)
$updates = @{
    2  = -4
    4  = -1
    6  = -2
    9  = -4
    26 = -2
    28 = -2
    29 = 0
    32 = -1
    38 = 6
    41 = -4
    43 = 0
    48 = -5
    50 = -3
    56 = -2
    60 = -3
    61 = -2
    63 = 6
    68 = -3
    69 = -1
    71 = -1
}

foreach ($row in $updates.Keys) {
    $ws.Range("F$row").Value = $updates[$row]
}
